$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.6003593248816514
$ws.Range("D2").Value = 0.1683500098719719
$ws.Range("E2").Value = 0.192509915913019
$ws.Range("F2").Value = 1.789372825043991
$ws.Range("G2").Value = 0.00249572892764612
$ws.Range("I2").Value = 1.370407801833476
$ws.Range("J2").Value = 0.2643003197562184
$ws.Range("K2").Value = 0.3851245660795257
$ws.Range("L2").Value = 0.2119011792620995
$ws.Range("M2").Value = 0.1808908174144683
$ws.Range("O2").Value = 4.491384077637321

# Row 3
$ws.Range("B3").Value = 0.5782418457902452
$ws.Range("D3").Value = 0.1682809498339921
$ws.Range("E3").Value = 0.1935080471101793
$ws.Range("F3").Value = 1.798807257976435
$ws.Range("G3").Value = 0.002498214679505963
$ws.Range("I3").Value = 1.383344258567494
$ws.Range("J3").Value = 0.2653397449350141
$ws.Range("K3").Value = 0.3375936572245735
$ws.Range("L3").Value = 0.2047353137760552
$ws.Range("M3").Value = 0.1744529858768153
$ws.Range("O3").Value = 4.518448774803474

# Row 4
$ws.Range("B4").Value = 0.5648414816183447
$ws.Range("D4").Value = 0.1682823996884864
$ws.Range("E4").Value = 0.1941590239807218
$ws.Range("F4").Value = 1.805394803104527
$ws.Range("G4").Value = 0.002499823765704405
$ws.Range("I4").Value = 1.391809298939613
$ws.Range("J4").Value = 0.2660175805884499
$ws.Range("K4").Value = 0.3083456780297524
$ws.Range("L4").Value = 0.2004047247810945
$ws.Range("M4").Value = 0.17055633264048
$ws.Range("O4").Value = 4.537178146587664

# Row 5
$ws.Range("B5").Value = 0.5594264888213445
$ws.Range("D5").Value = 0.1682940619846285
$ws.Range("E5").Value = 0.1944339089781262
$ws.Range("F5").Value = 1.808279331709549
$ws.Range("G5").Value = 0.002500500368986497
$ws.Range("I5").Value = 1.395390227981476
$ws.Range("J5").Value = 0.2663037934289321
$ws.Range("K5").Value = 0.2964115502476261
$ws.Range("L5").Value = 0.1986575459761468
$ws.Range("M5").Value = 0.1689826957761049
$ws.Range("O5").Value = 4.545341813303935

# Row 6
$ws.Range("B6").Value = 0.5585301132973086
$ws.Range("D6").Value = 0.1682966684924239
$ws.Range("E6").Value = 0.1944801342569846
$ws.Range("F6").Value = 1.808770393932797
$ws.Range("G6").Value = 0.002500613981741284
$ws.Range("I6").Value = 1.395992774895724
$ws.Range("J6").Value = 0.2663519229322429
$ws.Range("K6").Value = 0.2944289918381173
$ws.Range("L6").Value = 0.1983684939541348
$ws.Range("M6").Value = 0.1687222614895703
$ws.Range("O6").Value = 4.546729482671026

# Row 7
$ws.Range("B7").Value = 0.5647682671194332
$ws.Range("D7").Value = 0.1682825120841116
$ws.Range("E7").Value = 0.1941626922510062
$ws.Range("F7").Value = 1.805432894604053
$ws.Range("G7").Value = 0.002499832805784134
$ws.Range("I7").Value = 1.391857060618907
$ws.Range("J7").Value = 0.26602140007423
$ws.Range("K7").Value = 0.3081847913767035
$ws.Range("L7").Value = 0.2003810903713799
$ws.Range("M7").Value = 0.1705350520027871
$ws.Range("O7").Value = 4.537286093018267

# Row 8
$ws.Range("B8").Value = 0.5926962030878826
$ws.Range("D8").Value = 0.1683171225432432
$ws.Range("E8").Value = 0.1928461723427057
$ws.Range("F8").Value = 1.792460982863787
$ws.Range("G8").Value = 0.002496568861265465
$ws.Range("I8").Value = 1.374760007844301
$ws.Range("J8").Value = 0.2646505048201728
$ws.Range("K8").Value = 0.3687495612050213
$ws.Range("L8").Value = 0.2094160968279368
$ws.Range("M8").Value = 0.1786594753206678
$ws.Range("O8").Value = 4.500277931813258

# Row 9
$ws.Range("B9").Value = 0.6488693845241755
$ws.Range("D9").Value = 0.1687311121182304
$ws.Range("E9").Value = 0.1905659939757181
$ws.Range("F9").Value = 1.773321125483982
$ws.Range("G9").Value = 0.002490822671185614
$ws.Range("I9").Value = 1.345370130351959
$ws.Range("J9").Value = 0.2622754371421046
$ws.Range("K9").Value = 0.4869852402361516
$ws.Range("L9").Value = 0.2276778792362535
$ws.Range("M9").Value = 0.1950317611949508
$ws.Range("O9").Value = 4.444446286284887

# Row 10
$ws.Range("B10").Value = 0.6909743508040833
$ws.Range("D10").Value = 0.1692439677140101
$ws.Range("E10").Value = 0.1890732266980053
$ws.Range("F10").Value = 1.763088712950967
$ws.Range("G10").Value = 0.00248699598957092
$ws.Range("I10").Value = 1.326293150679504
$ws.Range("J10").Value = 0.2607198537269646
$ws.Range("K10").Value = 0.5735031825059025
$ws.Range("L10").Value = 0.2414206074014373
$ws.Range("M10").Value = 0.2073227087018026
$ws.Range("O10").Value = 4.413617968254869

# Row 11
$ws.Range("B11").Value = 0.7103056283417004
$ws.Range("D11").Value = 0.1695221189650269
$ws.Range("E11").Value = 0.1884334691165983
$ws.Range("F11").Value = 1.759263300650282
$ws.Range("G11").Value = 0.002485340085882635
$ws.Range("I11").Value = 1.318159428321071
$ws.Range("J11").Value = 0.2600529642494465
$ws.Range("K11").Value = 0.6127815528243445
$ws.Range("L11").Value = 0.2477421107184341
$ws.Range("M11").Value = 0.2129698907181208
$ws.Range("O11").Value = 4.401803389226416

# Row 12
$ws.Range("B12").Value = 0.7176509253268648
$ws.Range("D12").Value = 0.1696338547566825
$ws.Range("E12").Value = 0.1881968408286321
$ws.Range("F12").Value = 1.757933800694246
$ws.Range("G12").Value = 0.002484725180242314
$ws.Range("I12").Value = 1.315157603596383
$ws.Range("J12").Value = 0.2598062652004387
$ws.Range("K12").Value = 0.6276432762766433
$ws.Range("L12").Value = 0.250145815518465
$ws.Range("M12").Value = 0.2151162476190436
$ws.Range("O12").Value = 4.397646956355061

# Row 13
$ws.Range("B13").Value = 0.7160678829985727
$ws.Range("D13").Value = 0.1696095062401355
$ws.Range("E13").Value = 0.1882475526856844
$ws.Range("F13").Value = 1.75821483782687
$ws.Range("G13").Value = 0.002484857071582502
$ws.Range("I13").Value = 1.315800621456759
$ws.Range("J13").Value = 0.2598591369891796
$ws.Range("K13").Value = 0.6244430905510114
$ws.Range("L13").Value = 0.2496276973954679
$ws.Range("M13").Value = 0.2146536427639276
$ws.Range("O13").Value = 4.398528001422932

# Row 14
$ws.Range("B14").Value = 0.710909433305801
$ws.Range("D14").Value = 0.1695311834006574
$ws.Range("E14").Value = 0.1884138887572426
$ws.Range("F14").Value = 1.759151535855985
$ws.Range("G14").Value = 0.002485289254066026
$ws.Range("I14").Value = 1.317910898689256
$ws.Range("J14").Value = 0.2600325513068942
$ws.Range("K14").Value = 0.6140044835344156
$ws.Range("L14").Value = 0.2479396676796171
$ws.Range("M14").Value = 0.2131463155275952
$ws.Range("O14").Value = 4.401455075212084

# Row 15
$ws.Range("B15").Value = 0.7077529679331747
$ws.Range("D15").Value = 0.1694840413538827
$ws.Range("E15").Value = 0.1885165075679938
$ws.Range("F15").Value = 1.759740796090931
$ws.Range("G15").Value = 0.002485555558672232
$ws.Range("I15").Value = 1.319213691801224
$ws.Range("J15").Value = 0.2601395321301574
$ws.Range("K15").Value = 0.6076089310833481
$ws.Range("L15").Value = 0.2469069838482625
$ws.Range("M15").Value = 0.2122240571550691
$ws.Range("O15").Value = 4.403289332868212

# Row 16
$ws.Range("B16").Value = 0.6897145301179251
$ws.Range("D16").Value = 0.1692266886618015
$ws.Range("E16").Value = 0.1891158257744756
$ws.Range("F16").Value = 1.763355387814393
$ws.Range("G16").Value = 0.002487105909954301
$ws.Range("I16").Value = 1.32683566017116
$ws.Range("J16").Value = 0.2607642546325817
$ws.Range("K16").Value = 0.5709345913914206
$ws.Range("L16").Value = 0.2410088746772061
$ws.Range("M16").Value = 0.2069547649123251
$ws.Range("O16").Value = 4.41443451882256

# Row 17
$ws.Range("B17").Value = 0.6786935966760552
$ws.Range("D17").Value = 0.1690802674366765
$ws.Range("E17").Value = 0.1894935429282265
$ws.Range("F17").Value = 1.765785128301282
$ws.Range("G17").Value = 0.002488078698832114
$ws.Range("I17").Value = 1.331650906912458
$ws.Range("J17").Value = 0.2611579232827275
$ws.Range("K17").Value = 0.5484152687991468
$ws.Range("L17").Value = 0.2374083578320665
$ws.Range("M17").Value = 0.2037364481340447
$ws.Range("O17").Value = 4.421837451905645

# Row 18
$ws.Range("B18").Value = 0.6723713896896868
$ws.Range("D18").Value = 0.1690002739717826
$ws.Range("E18").Value = 0.1897144969732532
$ws.Range("F18").Value = 1.767260724005013
$ws.Range("G18").Value = 0.00248864621301271
$ws.Range("I18").Value = 1.334471763005446
$ws.Range("J18").Value = 0.2613881883110487
$ws.Range("K18").Value = 0.5354553531038277
$ws.Range("L18").Value = 0.2353440214550346
$ws.Range("M18").Value = 0.2018906347869702
$ws.Range("O18").Value = 4.426303393092979

# Row 19
$ws.Range("B19").Value = 0.6702336918060894
$ws.Range("D19").Value = 0.168973916389163
$ws.Range("E19").Value = 0.1897899444593083
$ws.Range("F19").Value = 1.76777374940837
$ws.Range("G19").Value = 0.002488839737994657
$ws.Range("I19").Value = 1.33543566102275
$ws.Range("J19").Value = 0.2614668118380177
$ws.Range("K19").Value = 0.5310661014535469
$ws.Range("L19").Value = 0.2346462085486962
$ws.Range("M19").Value = 0.2012665852674758
$ws.Range("O19").Value = 4.427851209589448

# Row 20
$ws.Range("B20").Value = 0.6798650651732316
$ws.Range("D20").Value = 0.1690954173581005
$ws.Range("E20").Value = 0.189452951375396
$ws.Range("F20").Value = 1.765518399154509
$ws.Range("G20").Value = 0.002487974317148483
$ws.Range("I20").Value = 1.331133011548573
$ws.Range("J20").Value = 0.261115619609301
$ws.Range("K20").Value = 0.5508132602883506
$ws.Range("L20").Value = 0.2377909585496809
$ws.Range("M20").Value = 0.2040784983992765
$ws.Range("O20").Value = 4.42102787448357

# Row 21
$ws.Range("B21").Value = 0.7124239220647155
$ws.Range("D21").Value = 0.1695540152075452
$ws.Range("E21").Value = 0.1883648790632919
$ws.Range("F21").Value = 1.758873173863883
$ws.Range("G21").Value = 0.002485161982256073
$ws.Range("I21").Value = 1.317288936175185
$ws.Range("J21").Value = 0.2599814570700065
$ws.Range("K21").Value = 0.6170708911485576
$ws.Range("L21").Value = 0.2484352158085414
$ws.Range("M21").Value = 0.2135888411444071
$ws.Range("O21").Value = 4.400586707759089

# Row 22
$ws.Range("B22").Value = 0.7338482521052754
$ws.Range("D22").Value = 0.1698910540955154
$ws.Range("E22").Value = 0.1876865917287485
$ws.Range("F22").Value = 1.755224270006707
$ws.Range("G22").Value = 0.002483394749091299
$ws.Range("I22").Value = 1.308697045310577
$ws.Range("J22").Value = 0.2592742314316521
$ws.Range("K22").Value = 0.6603029480127418
$ws.Range("L22").Value = 0.2554494188872809
$ws.Range("M22").Value = 0.2198503259565427
$ws.Range("O22").Value = 4.389077709101656

# Row 23
$ws.Range("B23").Value = 0.7224005812519749
$ws.Range("D23").Value = 0.169707769775691
$ws.Range("E23").Value = 0.1880456085597806
$ws.Range("F23").Value = 1.757108298453872
$ws.Range("G23").Value = 0.002484331495667453
$ws.Range("I23").Value = 1.313240992933952
$ws.Range("J23").Value = 0.2596485861359863
$ws.Range("K23").Value = 0.6372359490439123
$ws.Range("L23").Value = 0.251700591624072
$ws.Range("M23").Value = 0.2165043030722984
$ws.Range("O23").Value = 4.39505102959734

# Row 24
$ws.Range("B24").Value = 0.6793354010064263
$ws.Range("D24").Value = 0.1690885550369927
$ws.Range("E24").Value = 0.1894712909806291
$ws.Range("F24").Value = 1.76563874222709
$ws.Range("G24").Value = 0.002488021482285881
$ws.Range("I24").Value = 1.331366988483506
$ws.Range("J24").Value = 0.261134732827049
$ws.Range("K24").Value = 0.549729169573169
$ws.Range("L24").Value = 0.2376179671581156
$ws.Range("M24").Value = 0.2039238436331487
$ws.Range("O24").Value = 4.421393230603087

# Row 25
$ws.Range("B25").Value = 0.6335248461005847
$ws.Range("D25").Value = 0.1685823015339167
$ws.Range("E25").Value = 0.1911507015826102
$ws.Range("F25").Value = 1.777825730686715
$ws.Range("G25").Value = 0.002492307514715257
$ws.Range("I25").Value = 1.352878609103946
$ws.Range("J25").Value = 0.2628845849647838
$ws.Range("K25").Value = 0.455058938897082
$ws.Range("L25").Value = 0.2226799420802905
$ws.Range("M25").Value = 0.1905561274995051
$ws.Range("O25").Value = 4.457759481968651
